# This workbook lists daily Pimiento (pepper) price records (rows 2-88).
# The commit "Fruta / hortaliza, semanal" reorders/refreshes those data rows
# (same set of 87 records, rearranged) while leaving the header row (row 1)
# and the constant columns (A,B,C,E,F,G,R) untouched in content.
# We rebuild the A2:R88 block in one shot via a 2-D array write, which is
# both efficient and avoids any read/overwrite ordering issues.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 87,18
# target row 2 (was row 75)
$data[0,0] = 12
$data[0,1] = 'Mapocho Venta Directa de Santiago'
$data[0,2] = 'Metropolitana'
$data[0,3] = 44249
$data[0,4] = 13
$data[0,5] = 100112002
$data[0,6] = 'Pimiento'
$data[0,7] = 'Cuatro cascos verde'
$data[0,8] = 'Segunda'
$data[0,9] = 20
$data[0,10] = 6000
$data[0,11] = 6000
$data[0,12] = 6000
$data[0,13] = '$/caja 18 kilos'
$data[0,14] = 'Provincia de Limarí'
$data[0,15] = 333
$data[0,16] = 18
$data[0,17] = 'Hortaliza'
# target row 3 (was row 76)
$data[1,0] = 12
$data[1,1] = 'Mapocho Venta Directa de Santiago'
$data[1,2] = 'Metropolitana'
$data[1,3] = 44249
$data[1,4] = 13
$data[1,5] = 100112002
$data[1,6] = 'Pimiento'
$data[1,7] = 'Morrón rojo'
$data[1,8] = 'Primera'
$data[1,9] = 18
$data[1,10] = 15000
$data[1,11] = 15000
$data[1,12] = 15000
$data[1,13] = '$/caja 18 kilos'
$data[1,14] = 'Provincia de Limarí'
$data[1,15] = 833
$data[1,16] = 18
$data[1,17] = 'Hortaliza'
# target row 4 (was row 71)
$data[2,0] = 12
$data[2,1] = 'Mapocho Venta Directa de Santiago'
$data[2,2] = 'Metropolitana'
$data[2,3] = 44284
$data[2,4] = 13
$data[2,5] = 100112002
$data[2,6] = 'Pimiento'
$data[2,7] = 'Zafiro rojo'
$data[2,8] = 'Primera'
$data[2,9] = 20
$data[2,10] = 16000
$data[2,11] = 16000
$data[2,12] = 16000
$data[2,13] = '$/caja 18 kilos'
$data[2,14] = 'Provincia de Limarí'
$data[2,15] = 889
$data[2,16] = 18
$data[2,17] = 'Hortaliza'
# target row 5 (was row 72)
$data[3,0] = 12
$data[3,1] = 'Mapocho Venta Directa de Santiago'
$data[3,2] = 'Metropolitana'
$data[3,3] = 44284
$data[3,4] = 13
$data[3,5] = 100112002
$data[3,6] = 'Pimiento'
$data[3,7] = 'Zafiro verde'
$data[3,8] = 'Primera'
$data[3,9] = 25
$data[3,10] = 12000
$data[3,11] = 12000
$data[3,12] = 12000
$data[3,13] = '$/caja 18 kilos'
$data[3,14] = 'Provincia de Limarí'
$data[3,15] = 667
$data[3,16] = 18
$data[3,17] = 'Hortaliza'
# target row 6 (was row 73)
$data[4,0] = 12
$data[4,1] = 'Mapocho Venta Directa de Santiago'
$data[4,2] = 'Metropolitana'
$data[4,3] = 44277
$data[4,4] = 13
$data[4,5] = 100112002
$data[4,6] = 'Pimiento'
$data[4,7] = 'Zafiro rojo'
$data[4,8] = 'Primera'
$data[4,9] = 20
$data[4,10] = 16000
$data[4,11] = 16000
$data[4,12] = 16000
$data[4,13] = '$/caja 18 kilos'
$data[4,14] = 'Provincia de Limarí'
$data[4,15] = 889
$data[4,16] = 18
$data[4,17] = 'Hortaliza'
# target row 7 (was row 74)
$data[5,0] = 12
$data[5,1] = 'Mapocho Venta Directa de Santiago'
$data[5,2] = 'Metropolitana'
$data[5,3] = 44277
$data[5,4] = 13
$data[5,5] = 100112002
$data[5,6] = 'Pimiento'
$data[5,7] = 'Zafiro verde'
$data[5,8] = 'Primera'
$data[5,9] = 25
$data[5,10] = 12000
$data[5,11] = 12000
$data[5,12] = 12000
$data[5,13] = '$/caja 18 kilos'
$data[5,14] = 'Provincia de Limarí'
$data[5,15] = 667
$data[5,16] = 18
$data[5,17] = 'Hortaliza'
# target row 8 (was row 49)
$data[6,0] = 12
$data[6,1] = 'Mapocho Venta Directa de Santiago'
$data[6,2] = 'Metropolitana'
$data[6,3] = 44343
$data[6,4] = 13
$data[6,5] = 100112002
$data[6,6] = 'Pimiento'
$data[6,7] = 'Zafiro verde'
$data[6,8] = 'Primera'
$data[6,9] = 25
$data[6,10] = 14000
$data[6,11] = 14000
$data[6,12] = 14000
$data[6,13] = '$/caja 18 kilos'
$data[6,14] = 'Provincia de Limarí'
$data[6,15] = 778
$data[6,16] = 18
$data[6,17] = 'Hortaliza'
# target row 9 (was row 50)
$data[7,0] = 12
$data[7,1] = 'Mapocho Venta Directa de Santiago'
$data[7,2] = 'Metropolitana'
$data[7,3] = 44343
$data[7,4] = 13
$data[7,5] = 100112002
$data[7,6] = 'Pimiento'
$data[7,7] = 'Zafiro verde'
$data[7,8] = 'Segunda'
$data[7,9] = 15
$data[7,10] = 12000
$data[7,11] = 12000
$data[7,12] = 12000
$data[7,13] = '$/caja 18 kilos'
$data[7,14] = 'Provincia de Limarí'
$data[7,15] = 667
$data[7,16] = 18
$data[7,17] = 'Hortaliza'
# target row 10 (was row 43)
$data[8,0] = 12
$data[8,1] = 'Mapocho Venta Directa de Santiago'
$data[8,2] = 'Metropolitana'
$data[8,3] = 44425
$data[8,4] = 13
$data[8,5] = 100112002
$data[8,6] = 'Pimiento'
$data[8,7] = 'Morrón rojo'
$data[8,8] = 'Primera'
$data[8,9] = 8
$data[8,10] = 38000
$data[8,11] = 38000
$data[8,12] = 38000
$data[8,13] = '$/caja 18 kilos'
$data[8,14] = 'Provincia de Limarí'
$data[8,15] = 2111
$data[8,16] = 18
$data[8,17] = 'Hortaliza'
# target row 11 (was row 44)
$data[9,0] = 12
$data[9,1] = 'Mapocho Venta Directa de Santiago'
$data[9,2] = 'Metropolitana'
$data[9,3] = 44425
$data[9,4] = 13
$data[9,5] = 100112002
$data[9,6] = 'Pimiento'
$data[9,7] = 'Morrón rojo'
$data[9,8] = 'Segunda'
$data[9,9] = 15
$data[9,10] = 33000
$data[9,11] = 33000
$data[9,12] = 33000
$data[9,13] = '$/caja 18 kilos'
$data[9,14] = 'Provincia de Limarí'
$data[9,15] = 1833
$data[9,16] = 18
$data[9,17] = 'Hortaliza'
# target row 12 (was row 45)
$data[10,0] = 12
$data[10,1] = 'Mapocho Venta Directa de Santiago'
$data[10,2] = 'Metropolitana'
$data[10,3] = 44425
$data[10,4] = 13
$data[10,5] = 100112002
$data[10,6] = 'Pimiento'
$data[10,7] = 'Morrón rojo'
$data[10,8] = 'Tercera'
$data[10,9] = 10
$data[10,10] = 31000
$data[10,11] = 31000
$data[10,12] = 31000
$data[10,13] = '$/caja 18 kilos'
$data[10,14] = 'Provincia de Limarí'
$data[10,15] = 1722
$data[10,16] = 18
$data[10,17] = 'Hortaliza'
# target row 13 (was row 46)
$data[11,0] = 12
$data[11,1] = 'Mapocho Venta Directa de Santiago'
$data[11,2] = 'Metropolitana'
$data[11,3] = 44425
$data[11,4] = 13
$data[11,5] = 100112002
$data[11,6] = 'Pimiento'
$data[11,7] = 'Zafiro verde'
$data[11,8] = 'Primera'
$data[11,9] = 20
$data[11,10] = 35000
$data[11,11] = 35000
$data[11,12] = 35000
$data[11,13] = '$/caja 18 kilos'
$data[11,14] = 'Provincia de Limarí'
$data[11,15] = 1944
$data[11,16] = 18
$data[11,17] = 'Hortaliza'
# target row 14 (was row 86)
$data[12,0] = 12
$data[12,1] = 'Mapocho Venta Directa de Santiago'
$data[12,2] = 'Metropolitana'
$data[12,3] = 44474
$data[12,4] = 13
$data[12,5] = 100112002
$data[12,6] = 'Pimiento'
$data[12,7] = 'Cuatro cascos verde'
$data[12,8] = 'Primera'
$data[12,9] = 15
$data[12,10] = 38000
$data[12,11] = 38000
$data[12,12] = 38000
$data[12,13] = '$/caja 18 kilos'
$data[12,14] = 'Provincia de Limarí'
$data[12,15] = 2111
$data[12,16] = 18
$data[12,17] = 'Hortaliza'
# target row 15 (was row 87)
$data[13,0] = 12
$data[13,1] = 'Mapocho Venta Directa de Santiago'
$data[13,2] = 'Metropolitana'
$data[13,3] = 44474
$data[13,4] = 13
$data[13,5] = 100112002
$data[13,6] = 'Pimiento'
$data[13,7] = 'Cuatro cascos verde'
$data[13,8] = 'Segunda'
$data[13,9] = 25
$data[13,10] = 36000
$data[13,11] = 36000
$data[13,12] = 36000
$data[13,13] = '$/caja 18 kilos'
$data[13,14] = 'Provincia de Limarí'
$data[13,15] = 2000
$data[13,16] = 18
$data[13,17] = 'Hortaliza'
# target row 16 (was row 88)
$data[14,0] = 12
$data[14,1] = 'Mapocho Venta Directa de Santiago'
$data[14,2] = 'Metropolitana'
$data[14,3] = 44474
$data[14,4] = 13
$data[14,5] = 100112002
$data[14,6] = 'Pimiento'
$data[14,7] = 'Cuatro cascos verde'
$data[14,8] = 'Tercera'
$data[14,9] = 20
$data[14,10] = 34000
$data[14,11] = 34000
$data[14,12] = 34000
$data[14,13] = '$/caja 18 kilos'
$data[14,14] = 'Provincia de Limarí'
$data[14,15] = 1889
$data[14,16] = 18
$data[14,17] = 'Hortaliza'
# target row 17 (was row 10)
$data[15,0] = 12
$data[15,1] = 'Mapocho Venta Directa de Santiago'
$data[15,2] = 'Metropolitana'
$data[15,3] = 44193
$data[15,4] = 13
$data[15,5] = 100112002
$data[15,6] = 'Pimiento'
$data[15,7] = 'Zafiro rojo'
$data[15,8] = 'Tercera'
$data[15,9] = 20
$data[15,10] = 28000
$data[15,11] = 28000
$data[15,12] = 28000
$data[15,13] = '$/caja 18 kilos'
$data[15,14] = 'Provincia de Limarí'
$data[15,15] = 1556
$data[15,16] = 18
$data[15,17] = 'Hortaliza'
# target row 18 (was row 11)
$data[16,0] = 12
$data[16,1] = 'Mapocho Venta Directa de Santiago'
$data[16,2] = 'Metropolitana'
$data[16,3] = 44193
$data[16,4] = 13
$data[16,5] = 100112002
$data[16,6] = 'Pimiento'
$data[16,7] = 'Zafiro verde'
$data[16,8] = 'Primera'
$data[16,9] = 15
$data[16,10] = 18000
$data[16,11] = 18000
$data[16,12] = 18000
$data[16,13] = '$/caja 18 kilos'
$data[16,14] = 'Provincia de Limarí'
$data[16,15] = 1000
$data[16,16] = 18
$data[16,17] = 'Hortaliza'
# target row 19 (was row 12)
$data[17,0] = 12
$data[17,1] = 'Mapocho Venta Directa de Santiago'
$data[17,2] = 'Metropolitana'
$data[17,3] = 44193
$data[17,4] = 13
$data[17,5] = 100112002
$data[17,6] = 'Pimiento'
$data[17,7] = 'Zafiro verde'
$data[17,8] = 'Segunda'
$data[17,9] = 18
$data[17,10] = 16000
$data[17,11] = 16000
$data[17,12] = 16000
$data[17,13] = '$/caja 18 kilos'
$data[17,14] = 'Provincia de Limarí'
$data[17,15] = 889
$data[17,16] = 18
$data[17,17] = 'Hortaliza'
# target row 20 (was row 77)
$data[18,0] = 12
$data[18,1] = 'Mapocho Venta Directa de Santiago'
$data[18,2] = 'Metropolitana'
$data[18,3] = 44421
$data[18,4] = 13
$data[18,5] = 100112002
$data[18,6] = 'Pimiento'
$data[18,7] = 'Zafiro rojo'
$data[18,8] = 'Primera'
$data[18,9] = 15
$data[18,10] = 28000
$data[18,11] = 28000
$data[18,12] = 28000
$data[18,13] = '$/caja 18 kilos'
$data[18,14] = 'Provincia de Limarí'
$data[18,15] = 1556
$data[18,16] = 18
$data[18,17] = 'Hortaliza'
# target row 21 (was row 78)
$data[19,0] = 12
$data[19,1] = 'Mapocho Venta Directa de Santiago'
$data[19,2] = 'Metropolitana'
$data[19,3] = 44421
$data[19,4] = 13
$data[19,5] = 100112002
$data[19,6] = 'Pimiento'
$data[19,7] = 'Zafiro rojo'
$data[19,8] = 'Segunda'
$data[19,9] = 20
$data[19,10] = 26000
$data[19,11] = 26000
$data[19,12] = 26000
$data[19,13] = '$/caja 18 kilos'
$data[19,14] = 'Provincia de Limarí'
$data[19,15] = 1444
$data[19,16] = 18
$data[19,17] = 'Hortaliza'
# target row 22 (was row 79)
$data[20,0] = 12
$data[20,1] = 'Mapocho Venta Directa de Santiago'
$data[20,2] = 'Metropolitana'
$data[20,3] = 44421
$data[20,4] = 13
$data[20,5] = 100112002
$data[20,6] = 'Pimiento'
$data[20,7] = 'Zafiro verde'
$data[20,8] = 'Primera'
$data[20,9] = 15
$data[20,10] = 32000
$data[20,11] = 32000
$data[20,12] = 32000
$data[20,13] = '$/caja 18 kilos'
$data[20,14] = 'Provincia de Limarí'
$data[20,15] = 1778
$data[20,16] = 18
$data[20,17] = 'Hortaliza'
# target row 23 (was row 80)
$data[21,0] = 12
$data[21,1] = 'Mapocho Venta Directa de Santiago'
$data[21,2] = 'Metropolitana'
$data[21,3] = 44421
$data[21,4] = 13
$data[21,5] = 100112002
$data[21,6] = 'Pimiento'
$data[21,7] = 'Zafiro verde'
$data[21,8] = 'Segunda'
$data[21,9] = 10
$data[21,10] = 30000
$data[21,11] = 30000
$data[21,12] = 30000
$data[21,13] = '$/caja 18 kilos'
$data[21,14] = 'Provincia de Limarí'
$data[21,15] = 1667
$data[21,16] = 18
$data[21,17] = 'Hortaliza'
# target row 24 (was row 81)
$data[22,0] = 12
$data[22,1] = 'Mapocho Venta Directa de Santiago'
$data[22,2] = 'Metropolitana'
$data[22,3] = 44421
$data[22,4] = 13
$data[22,5] = 100112002
$data[22,6] = 'Pimiento'
$data[22,7] = 'Zafiro verde'
$data[22,8] = 'Tercera'
$data[22,9] = 12
$data[22,10] = 28000
$data[22,11] = 28000
$data[22,12] = 28000
$data[22,13] = '$/caja 18 kilos'
$data[22,14] = 'Provincia de Limarí'
$data[22,15] = 1556
$data[22,16] = 18
$data[22,17] = 'Hortaliza'
# target row 25 (was row 82)
$data[23,0] = 12
$data[23,1] = 'Mapocho Venta Directa de Santiago'
$data[23,2] = 'Metropolitana'
$data[23,3] = 44291
$data[23,4] = 13
$data[23,5] = 100112002
$data[23,6] = 'Pimiento'
$data[23,7] = 'Morrón rojo'
$data[23,8] = 'Primera'
$data[23,9] = 20
$data[23,10] = 10000
$data[23,11] = 10000
$data[23,12] = 10000
$data[23,13] = '$/caja 18 kilos'
$data[23,14] = 'Provincia de Limarí'
$data[23,15] = 556
$data[23,16] = 18
$data[23,17] = 'Hortaliza'
# target row 26 (was row 21)
$data[24,0] = 12
$data[24,1] = 'Mapocho Venta Directa de Santiago'
$data[24,2] = 'Metropolitana'
$data[24,3] = 44460
$data[24,4] = 13
$data[24,5] = 100112002
$data[24,6] = 'Pimiento'
$data[24,7] = 'Zafiro verde'
$data[24,8] = 'Primera'
$data[24,9] = 55
$data[24,10] = 36000
$data[24,11] = 36000
$data[24,12] = 36000
$data[24,13] = '$/caja 18 kilos'
$data[24,14] = 'Provincia de Limarí'
$data[24,15] = 2000
$data[24,16] = 18
$data[24,17] = 'Hortaliza'
# target row 27 (was row 39)
$data[25,0] = 12
$data[25,1] = 'Mapocho Venta Directa de Santiago'
$data[25,2] = 'Metropolitana'
$data[25,3] = 44280
$data[25,4] = 13
$data[25,5] = 100112002
$data[25,6] = 'Pimiento'
$data[25,7] = 'Zafiro rojo'
$data[25,8] = 'Primera'
$data[25,9] = 30
$data[25,10] = 16000
$data[25,11] = 16000
$data[25,12] = 16000
$data[25,13] = '$/caja 18 kilos'
$data[25,14] = 'Provincia de Limarí'
$data[25,15] = 889
$data[25,16] = 18
$data[25,17] = 'Hortaliza'
# target row 28 (was row 40)
$data[26,0] = 12
$data[26,1] = 'Mapocho Venta Directa de Santiago'
$data[26,2] = 'Metropolitana'
$data[26,3] = 44280
$data[26,4] = 13
$data[26,5] = 100112002
$data[26,6] = 'Pimiento'
$data[26,7] = 'Zafiro verde'
$data[26,8] = 'Primera'
$data[26,9] = 45
$data[26,10] = 12000
$data[26,11] = 12000
$data[26,12] = 12000
$data[26,13] = '$/caja 18 kilos'
$data[26,14] = 'Provincia de Limarí'
$data[26,15] = 667
$data[26,16] = 18
$data[26,17] = 'Hortaliza'
# target row 29 (was row 30)
$data[27,0] = 12
$data[27,1] = 'Mapocho Venta Directa de Santiago'
$data[27,2] = 'Metropolitana'
$data[27,3] = 44270
$data[27,4] = 13
$data[27,5] = 100112002
$data[27,6] = 'Pimiento'
$data[27,7] = 'Zafiro rojo'
$data[27,8] = 'Primera'
$data[27,9] = 25
$data[27,10] = 17000
$data[27,11] = 17000
$data[27,12] = 17000
$data[27,13] = '$/caja 18 kilos'
$data[27,14] = 'Provincia de Limarí'
$data[27,15] = 944
$data[27,16] = 18
$data[27,17] = 'Hortaliza'
# target row 30 (was row 31)
$data[28,0] = 12
$data[28,1] = 'Mapocho Venta Directa de Santiago'
$data[28,2] = 'Metropolitana'
$data[28,3] = 44270
$data[28,4] = 13
$data[28,5] = 100112002
$data[28,6] = 'Pimiento'
$data[28,7] = 'Zafiro verde'
$data[28,8] = 'Primera'
$data[28,9] = 30
$data[28,10] = 10000
$data[28,11] = 10000
$data[28,12] = 10000
$data[28,13] = '$/caja 18 kilos'
$data[28,14] = 'Provincia de Limarí'
$data[28,15] = 556
$data[28,16] = 18
$data[28,17] = 'Hortaliza'
# target row 31 (was row 34)
$data[29,0] = 12
$data[29,1] = 'Mapocho Venta Directa de Santiago'
$data[29,2] = 'Metropolitana'
$data[29,3] = 44186
$data[29,4] = 13
$data[29,5] = 100112002
$data[29,6] = 'Pimiento'
$data[29,7] = 'Zafiro verde'
$data[29,8] = 'Primera'
$data[29,9] = 20
$data[29,10] = 17000
$data[29,11] = 17000
$data[29,12] = 17000
$data[29,13] = '$/caja 18 kilos'
$data[29,14] = 'Provincia de Limarí'
$data[29,15] = 944
$data[29,16] = 18
$data[29,17] = 'Hortaliza'
# target row 32 (was row 54)
$data[30,0] = 12
$data[30,1] = 'Mapocho Venta Directa de Santiago'
$data[30,2] = 'Metropolitana'
$data[30,3] = 44432
$data[30,4] = 13
$data[30,5] = 100112002
$data[30,6] = 'Pimiento'
$data[30,7] = 'Zafiro verde'
$data[30,8] = 'Primera'
$data[30,9] = 20
$data[30,10] = 35000
$data[30,11] = 35000
$data[30,12] = 35000
$data[30,13] = '$/caja 18 kilos'
$data[30,14] = 'Provincia de Limarí'
$data[30,15] = 1944
$data[30,16] = 18
$data[30,17] = 'Hortaliza'
# target row 33 (was row 13)
$data[31,0] = 12
$data[31,1] = 'Mapocho Venta Directa de Santiago'
$data[31,2] = 'Metropolitana'
$data[31,3] = 44312
$data[31,4] = 13
$data[31,5] = 100112002
$data[31,6] = 'Pimiento'
$data[31,7] = 'Zafiro rojo'
$data[31,8] = 'Primera'
$data[31,9] = 25
$data[31,10] = 20000
$data[31,11] = 20000
$data[31,12] = 20000
$data[31,13] = '$/caja 18 kilos'
$data[31,14] = 'Provincia de Limarí'
$data[31,15] = 1111
$data[31,16] = 18
$data[31,17] = 'Hortaliza'
# target row 34 (was row 14)
$data[32,0] = 12
$data[32,1] = 'Mapocho Venta Directa de Santiago'
$data[32,2] = 'Metropolitana'
$data[32,3] = 44312
$data[32,4] = 13
$data[32,5] = 100112002
$data[32,6] = 'Pimiento'
$data[32,7] = 'Zafiro verde'
$data[32,8] = 'Primera'
$data[32,9] = 30
$data[32,10] = 15000
$data[32,11] = 15000
$data[32,12] = 15000
$data[32,13] = '$/caja 18 kilos'
$data[32,14] = 'Provincia de Limarí'
$data[32,15] = 833
$data[32,16] = 18
$data[32,17] = 'Hortaliza'
# target row 35 (was row 22)
$data[33,0] = 12
$data[33,1] = 'Mapocho Venta Directa de Santiago'
$data[33,2] = 'Metropolitana'
$data[33,3] = 44581
$data[33,4] = 13
$data[33,5] = 100112002
$data[33,6] = 'Pimiento'
$data[33,7] = 'Zafiro rojo'
$data[33,8] = 'Primera'
$data[33,9] = 20
$data[33,10] = 18000
$data[33,11] = 18000
$data[33,12] = 18000
$data[33,13] = '$/caja 18 kilos'
$data[33,14] = 'Provincia de Limarí'
$data[33,15] = 1000
$data[33,16] = 18
$data[33,17] = 'Hortaliza'
# target row 36 (was row 23)
$data[34,0] = 12
$data[34,1] = 'Mapocho Venta Directa de Santiago'
$data[34,2] = 'Metropolitana'
$data[34,3] = 44581
$data[34,4] = 13
$data[34,5] = 100112002
$data[34,6] = 'Pimiento'
$data[34,7] = 'Zafiro rojo'
$data[34,8] = 'Segunda'
$data[34,9] = 20
$data[34,10] = 15000
$data[34,11] = 15000
$data[34,12] = 15000
$data[34,13] = '$/caja 18 kilos'
$data[34,14] = 'Provincia de Limarí'
$data[34,15] = 833
$data[34,16] = 18
$data[34,17] = 'Hortaliza'
# target row 37 (was row 24)
$data[35,0] = 12
$data[35,1] = 'Mapocho Venta Directa de Santiago'
$data[35,2] = 'Metropolitana'
$data[35,3] = 44581
$data[35,4] = 13
$data[35,5] = 100112002
$data[35,6] = 'Pimiento'
$data[35,7] = 'Zafiro verde'
$data[35,8] = 'Primera'
$data[35,9] = 20
$data[35,10] = 14000
$data[35,11] = 14000
$data[35,12] = 14000
$data[35,13] = '$/caja 18 kilos'
$data[35,14] = 'Provincia de Limarí'
$data[35,15] = 778
$data[35,16] = 18
$data[35,17] = 'Hortaliza'
# target row 38 (was row 25)
$data[36,0] = 12
$data[36,1] = 'Mapocho Venta Directa de Santiago'
$data[36,2] = 'Metropolitana'
$data[36,3] = 44581
$data[36,4] = 13
$data[36,5] = 100112002
$data[36,6] = 'Pimiento'
$data[36,7] = 'Zafiro verde'
$data[36,8] = 'Segunda'
$data[36,9] = 35
$data[36,10] = 12000
$data[36,11] = 12000
$data[36,12] = 12000
$data[36,13] = '$/caja 18 kilos'
$data[36,14] = 'Provincia de Limarí'
$data[36,15] = 667
$data[36,16] = 18
$data[36,17] = 'Hortaliza'
# target row 39 (was row 33)
$data[37,0] = 12
$data[37,1] = 'Mapocho Venta Directa de Santiago'
$data[37,2] = 'Metropolitana'
$data[37,3] = 44333
$data[37,4] = 13
$data[37,5] = 100112002
$data[37,6] = 'Pimiento'
$data[37,7] = 'Zafiro verde'
$data[37,8] = 'Primera'
$data[37,9] = 25
$data[37,10] = 12000
$data[37,11] = 13000
$data[37,12] = 12600
$data[37,13] = '$/caja 18 kilos'
$data[37,14] = 'Provincia de Limarí'
$data[37,15] = 700
$data[37,16] = 18
$data[37,17] = 'Hortaliza'
# target row 40 (was row 32)
$data[38,0] = 12
$data[38,1] = 'Mapocho Venta Directa de Santiago'
$data[38,2] = 'Metropolitana'
$data[38,3] = 44376
$data[38,4] = 13
$data[38,5] = 100112002
$data[38,6] = 'Pimiento'
$data[38,7] = 'Zafiro verde'
$data[38,8] = 'Primera'
$data[38,9] = 30
$data[38,10] = 16000
$data[38,11] = 16000
$data[38,12] = 16000
$data[38,13] = '$/caja 18 kilos'
$data[38,14] = 'Provincia de Limarí'
$data[38,15] = 889
$data[38,16] = 18
$data[38,17] = 'Hortaliza'
# target row 41 (was row 3)
$data[39,0] = 12
$data[39,1] = 'Mapocho Venta Directa de Santiago'
$data[39,2] = 'Metropolitana'
$data[39,3] = 44236
$data[39,4] = 13
$data[39,5] = 100112002
$data[39,6] = 'Pimiento'
$data[39,7] = 'Cuatro cascos rojo'
$data[39,8] = 'Extra'
$data[39,9] = 60
$data[39,10] = 25000
$data[39,11] = 25000
$data[39,12] = 25000
$data[39,13] = '$/caja 18 kilos'
$data[39,14] = 'Provincia de Limarí'
$data[39,15] = 1389
$data[39,16] = 18
$data[39,17] = 'Hortaliza'
# target row 42 (was row 4)
$data[40,0] = 12
$data[40,1] = 'Mapocho Venta Directa de Santiago'
$data[40,2] = 'Metropolitana'
$data[40,3] = 44236
$data[40,4] = 13
$data[40,5] = 100112002
$data[40,6] = 'Pimiento'
$data[40,7] = 'Cuatro cascos rojo'
$data[40,8] = 'Primera'
$data[40,9] = 120
$data[40,10] = 23000
$data[40,11] = 23000
$data[40,12] = 23000
$data[40,13] = '$/caja 18 kilos'
$data[40,14] = 'Provincia de Limarí'
$data[40,15] = 1278
$data[40,16] = 18
$data[40,17] = 'Hortaliza'
# target row 43 (was row 5)
$data[41,0] = 12
$data[41,1] = 'Mapocho Venta Directa de Santiago'
$data[41,2] = 'Metropolitana'
$data[41,3] = 44236
$data[41,4] = 13
$data[41,5] = 100112002
$data[41,6] = 'Pimiento'
$data[41,7] = 'Cuatro cascos rojo'
$data[41,8] = 'Segunda'
$data[41,9] = 80
$data[41,10] = 21000
$data[41,11] = 21000
$data[41,12] = 21000
$data[41,13] = '$/caja 18 kilos'
$data[41,14] = 'Provincia de Limarí'
$data[41,15] = 1167
$data[41,16] = 18
$data[41,17] = 'Hortaliza'
# target row 44 (was row 6)
$data[42,0] = 12
$data[42,1] = 'Mapocho Venta Directa de Santiago'
$data[42,2] = 'Metropolitana'
$data[42,3] = 44236
$data[42,4] = 13
$data[42,5] = 100112002
$data[42,6] = 'Pimiento'
$data[42,7] = 'Cuatro cascos rojo'
$data[42,8] = 'Tercera'
$data[42,9] = 50
$data[42,10] = 18000
$data[42,11] = 18000
$data[42,12] = 18000
$data[42,13] = '$/caja 18 kilos'
$data[42,14] = 'Provincia de Limarí'
$data[42,15] = 1000
$data[42,16] = 18
$data[42,17] = 'Hortaliza'
# target row 45 (was row 7)
$data[43,0] = 12
$data[43,1] = 'Mapocho Venta Directa de Santiago'
$data[43,2] = 'Metropolitana'
$data[43,3] = 44236
$data[43,4] = 13
$data[43,5] = 100112002
$data[43,6] = 'Pimiento'
$data[43,7] = 'Cuatro cascos verde'
$data[43,8] = 'Primera'
$data[43,9] = 150
$data[43,10] = 10000
$data[43,11] = 10000
$data[43,12] = 10000
$data[43,13] = '$/caja 18 kilos'
$data[43,14] = 'Provincia de Limarí'
$data[43,15] = 556
$data[43,16] = 18
$data[43,17] = 'Hortaliza'
# target row 46 (was row 8)
$data[44,0] = 12
$data[44,1] = 'Mapocho Venta Directa de Santiago'
$data[44,2] = 'Metropolitana'
$data[44,3] = 44236
$data[44,4] = 13
$data[44,5] = 100112002
$data[44,6] = 'Pimiento'
$data[44,7] = 'Cuatro cascos verde'
$data[44,8] = 'Segunda'
$data[44,9] = 100
$data[44,10] = 8000
$data[44,11] = 8000
$data[44,12] = 8000
$data[44,13] = '$/caja 18 kilos'
$data[44,14] = 'Provincia de Limarí'
$data[44,15] = 444
$data[44,16] = 18
$data[44,17] = 'Hortaliza'
# target row 47 (was row 9)
$data[45,0] = 12
$data[45,1] = 'Mapocho Venta Directa de Santiago'
$data[45,2] = 'Metropolitana'
$data[45,3] = 44236
$data[45,4] = 13
$data[45,5] = 100112002
$data[45,6] = 'Pimiento'
$data[45,7] = 'Cuatro cascos verde'
$data[45,8] = 'Tercera'
$data[45,9] = 75
$data[45,10] = 6000
$data[45,11] = 6000
$data[45,12] = 6000
$data[45,13] = '$/caja 18 kilos'
$data[45,14] = 'Provincia de Limarí'
$data[45,15] = 333
$data[45,16] = 18
$data[45,17] = 'Hortaliza'
# target row 48 (was row 35)
$data[46,0] = 12
$data[46,1] = 'Mapocho Venta Directa de Santiago'
$data[46,2] = 'Metropolitana'
$data[46,3] = 44446
$data[46,4] = 13
$data[46,5] = 100112002
$data[46,6] = 'Pimiento'
$data[46,7] = 'Zafiro rojo'
$data[46,8] = 'Primera'
$data[46,9] = 10
$data[46,10] = 34000
$data[46,11] = 34000
$data[46,12] = 34000
$data[46,13] = '$/caja 18 kilos'
$data[46,14] = 'Provincia de Limarí'
$data[46,15] = 1889
$data[46,16] = 18
$data[46,17] = 'Hortaliza'
# target row 49 (was row 36)
$data[47,0] = 12
$data[47,1] = 'Mapocho Venta Directa de Santiago'
$data[47,2] = 'Metropolitana'
$data[47,3] = 44446
$data[47,4] = 13
$data[47,5] = 100112002
$data[47,6] = 'Pimiento'
$data[47,7] = 'Zafiro verde'
$data[47,8] = 'Primera'
$data[47,9] = 10
$data[47,10] = 33000
$data[47,11] = 33000
$data[47,12] = 33000
$data[47,13] = '$/caja 18 kilos'
$data[47,14] = 'Provincia de Limarí'
$data[47,15] = 1833
$data[47,16] = 18
$data[47,17] = 'Hortaliza'
# target row 50 (was row 37)
$data[48,0] = 12
$data[48,1] = 'Mapocho Venta Directa de Santiago'
$data[48,2] = 'Metropolitana'
$data[48,3] = 44446
$data[48,4] = 13
$data[48,5] = 100112002
$data[48,6] = 'Pimiento'
$data[48,7] = 'Zafiro verde'
$data[48,8] = 'Segunda'
$data[48,9] = 8
$data[48,10] = 31000
$data[48,11] = 31000
$data[48,12] = 31000
$data[48,13] = '$/caja 18 kilos'
$data[48,14] = 'Provincia de Limarí'
$data[48,15] = 1722
$data[48,16] = 18
$data[48,17] = 'Hortaliza'
# target row 51 (was row 38)
$data[49,0] = 12
$data[49,1] = 'Mapocho Venta Directa de Santiago'
$data[49,2] = 'Metropolitana'
$data[49,3] = 44446
$data[49,4] = 13
$data[49,5] = 100112002
$data[49,6] = 'Pimiento'
$data[49,7] = 'Zafiro verde'
$data[49,8] = 'Tercera'
$data[49,9] = 12
$data[49,10] = 29000
$data[49,11] = 29000
$data[49,12] = 29000
$data[49,13] = '$/caja 18 kilos'
$data[49,14] = 'Provincia de Limarí'
$data[49,15] = 1611
$data[49,16] = 18
$data[49,17] = 'Hortaliza'
# target row 52 (was row 2)
$data[50,0] = 12
$data[50,1] = 'Mapocho Venta Directa de Santiago'
$data[50,2] = 'Metropolitana'
$data[50,3] = 44232
$data[50,4] = 13
$data[50,5] = 100112002
$data[50,6] = 'Pimiento'
$data[50,7] = 'Cuatro cascos verde'
$data[50,8] = 'Primera'
$data[50,9] = 70
$data[50,10] = 12000
$data[50,11] = 12000
$data[50,12] = 12000
$data[50,13] = '$/caja 15 kilos'
$data[50,14] = 'Región de Arica y Parinacota'
$data[50,15] = 800
$data[50,16] = 15
$data[50,17] = 'Hortaliza'
# target row 53 (was row 47)
$data[51,0] = 12
$data[51,1] = 'Mapocho Venta Directa de Santiago'
$data[51,2] = 'Metropolitana'
$data[51,3] = 44585
$data[51,4] = 13
$data[51,5] = 100112002
$data[51,6] = 'Pimiento'
$data[51,7] = 'Zafiro rojo'
$data[51,8] = 'Primera'
$data[51,9] = 25
$data[51,10] = 17000
$data[51,11] = 17000
$data[51,12] = 17000
$data[51,13] = '$/caja 18 kilos'
$data[51,14] = 'Provincia de Limarí'
$data[51,15] = 944
$data[51,16] = 18
$data[51,17] = 'Hortaliza'
# target row 54 (was row 48)
$data[52,0] = 12
$data[52,1] = 'Mapocho Venta Directa de Santiago'
$data[52,2] = 'Metropolitana'
$data[52,3] = 44585
$data[52,4] = 13
$data[52,5] = 100112002
$data[52,6] = 'Pimiento'
$data[52,7] = 'Zafiro rojo'
$data[52,8] = 'Segunda'
$data[52,9] = 30
$data[52,10] = 14000
$data[52,11] = 14000
$data[52,12] = 14000
$data[52,13] = '$/caja 18 kilos'
$data[52,14] = 'Provincia de Limarí'
$data[52,15] = 778
$data[52,16] = 18
$data[52,17] = 'Hortaliza'
# target row 55 (was row 28)
$data[53,0] = 12
$data[53,1] = 'Mapocho Venta Directa de Santiago'
$data[53,2] = 'Metropolitana'
$data[53,3] = 44243
$data[53,4] = 13
$data[53,5] = 100112002
$data[53,6] = 'Pimiento'
$data[53,7] = 'Cuatro cascos rojo'
$data[53,8] = 'Primera'
$data[53,9] = 55
$data[53,10] = 20000
$data[53,11] = 22000
$data[53,12] = 21091
$data[53,13] = '$/caja 18 kilos'
$data[53,14] = 'Provincia de Quillota'
$data[53,15] = 1172
$data[53,16] = 18
$data[53,17] = 'Hortaliza'
# target row 56 (was row 29)
$data[54,0] = 12
$data[54,1] = 'Mapocho Venta Directa de Santiago'
$data[54,2] = 'Metropolitana'
$data[54,3] = 44243
$data[54,4] = 13
$data[54,5] = 100112002
$data[54,6] = 'Pimiento'
$data[54,7] = 'Cuatro cascos verde'
$data[54,8] = 'Primera'
$data[54,9] = 90
$data[54,10] = 12000
$data[54,11] = 13000
$data[54,12] = 12556
$data[54,13] = '$/caja 18 kilos'
$data[54,14] = 'Provincia de Quillota'
$data[54,15] = 698
$data[54,16] = 18
$data[54,17] = 'Hortaliza'
# target row 57 (was row 61)
$data[55,0] = 12
$data[55,1] = 'Mapocho Venta Directa de Santiago'
$data[55,2] = 'Metropolitana'
$data[55,3] = 44435
$data[55,4] = 13
$data[55,5] = 100112002
$data[55,6] = 'Pimiento'
$data[55,7] = 'Zafiro verde'
$data[55,8] = 'Primera'
$data[55,9] = 30
$data[55,10] = 35000
$data[55,11] = 35000
$data[55,12] = 35000
$data[55,13] = '$/caja 18 kilos'
$data[55,14] = 'Provincia de Limarí'
$data[55,15] = 1944
$data[55,16] = 18
$data[55,17] = 'Hortaliza'
# target row 58 (was row 66)
$data[56,0] = 12
$data[56,1] = 'Mapocho Venta Directa de Santiago'
$data[56,2] = 'Metropolitana'
$data[56,3] = 44418
$data[56,4] = 13
$data[56,5] = 100112002
$data[56,6] = 'Pimiento'
$data[56,7] = 'Morrón rojo'
$data[56,8] = 'Primera'
$data[56,9] = 10
$data[56,10] = 33000
$data[56,11] = 33000
$data[56,12] = 33000
$data[56,13] = '$/caja 18 kilos'
$data[56,14] = 'Provincia de Limarí'
$data[56,15] = 1833
$data[56,16] = 18
$data[56,17] = 'Hortaliza'
# target row 59 (was row 67)
$data[57,0] = 12
$data[57,1] = 'Mapocho Venta Directa de Santiago'
$data[57,2] = 'Metropolitana'
$data[57,3] = 44418
$data[57,4] = 13
$data[57,5] = 100112002
$data[57,6] = 'Pimiento'
$data[57,7] = 'Morrón rojo'
$data[57,8] = 'Segunda'
$data[57,9] = 12
$data[57,10] = 31000
$data[57,11] = 31000
$data[57,12] = 31000
$data[57,13] = '$/caja 18 kilos'
$data[57,14] = 'Provincia de Limarí'
$data[57,15] = 1722
$data[57,16] = 18
$data[57,17] = 'Hortaliza'
# target row 60 (was row 68)
$data[58,0] = 12
$data[58,1] = 'Mapocho Venta Directa de Santiago'
$data[58,2] = 'Metropolitana'
$data[58,3] = 44418
$data[58,4] = 13
$data[58,5] = 100112002
$data[58,6] = 'Pimiento'
$data[58,7] = 'Zafiro rojo'
$data[58,8] = 'Primera'
$data[58,9] = 10
$data[58,10] = 28000
$data[58,11] = 28000
$data[58,12] = 28000
$data[58,13] = '$/caja 18 kilos'
$data[58,14] = 'Provincia de Limarí'
$data[58,15] = 1556
$data[58,16] = 18
$data[58,17] = 'Hortaliza'
# target row 61 (was row 69)
$data[59,0] = 12
$data[59,1] = 'Mapocho Venta Directa de Santiago'
$data[59,2] = 'Metropolitana'
$data[59,3] = 44418
$data[59,4] = 13
$data[59,5] = 100112002
$data[59,6] = 'Pimiento'
$data[59,7] = 'Zafiro rojo'
$data[59,8] = 'Segunda'
$data[59,9] = 15
$data[59,10] = 26000
$data[59,11] = 26000
$data[59,12] = 26000
$data[59,13] = '$/caja 18 kilos'
$data[59,14] = 'Provincia de Limarí'
$data[59,15] = 1444
$data[59,16] = 18
$data[59,17] = 'Hortaliza'
# target row 62 (was row 70)
$data[60,0] = 12
$data[60,1] = 'Mapocho Venta Directa de Santiago'
$data[60,2] = 'Metropolitana'
$data[60,3] = 44418
$data[60,4] = 13
$data[60,5] = 100112002
$data[60,6] = 'Pimiento'
$data[60,7] = 'Zafiro verde'
$data[60,8] = 'Primera'
$data[60,9] = 10
$data[60,10] = 28000
$data[60,11] = 28000
$data[60,12] = 28000
$data[60,13] = '$/caja 18 kilos'
$data[60,14] = 'Provincia de Limarí'
$data[60,15] = 1556
$data[60,16] = 18
$data[60,17] = 'Hortaliza'
# target row 63 (was row 15)
$data[61,0] = 12
$data[61,1] = 'Mapocho Venta Directa de Santiago'
$data[61,2] = 'Metropolitana'
$data[61,3] = 44179
$data[61,4] = 13
$data[61,5] = 100112002
$data[61,6] = 'Pimiento'
$data[61,7] = 'Zafiro rojo'
$data[61,8] = 'Tercera'
$data[61,9] = 20
$data[61,10] = 22000
$data[61,11] = 22000
$data[61,12] = 22000
$data[61,13] = '$/caja 18 kilos'
$data[61,14] = 'Provincia de Limarí'
$data[61,15] = 1222
$data[61,16] = 18
$data[61,17] = 'Hortaliza'
# target row 64 (was row 16)
$data[62,0] = 12
$data[62,1] = 'Mapocho Venta Directa de Santiago'
$data[62,2] = 'Metropolitana'
$data[62,3] = 44179
$data[62,4] = 13
$data[62,5] = 100112002
$data[62,6] = 'Pimiento'
$data[62,7] = 'Zafiro verde'
$data[62,8] = 'Primera'
$data[62,9] = 25
$data[62,10] = 16000
$data[62,11] = 16000
$data[62,12] = 16000
$data[62,13] = '$/caja 18 kilos'
$data[62,14] = 'Provincia de Limarí'
$data[62,15] = 889
$data[62,16] = 18
$data[62,17] = 'Hortaliza'
# target row 65 (was row 17)
$data[63,0] = 12
$data[63,1] = 'Mapocho Venta Directa de Santiago'
$data[63,2] = 'Metropolitana'
$data[63,3] = 44179
$data[63,4] = 13
$data[63,5] = 100112002
$data[63,6] = 'Pimiento'
$data[63,7] = 'Zafiro verde'
$data[63,8] = 'Segunda'
$data[63,9] = 20
$data[63,10] = 14000
$data[63,11] = 14000
$data[63,12] = 14000
$data[63,13] = '$/caja 18 kilos'
$data[63,14] = 'Provincia de Limarí'
$data[63,15] = 778
$data[63,16] = 18
$data[63,17] = 'Hortaliza'
# target row 66 (was row 57)
$data[64,0] = 12
$data[64,1] = 'Mapocho Venta Directa de Santiago'
$data[64,2] = 'Metropolitana'
$data[64,3] = 44467
$data[64,4] = 13
$data[64,5] = 100112002
$data[64,6] = 'Pimiento'
$data[64,7] = 'Cuatro cascos verde'
$data[64,8] = 'Primera'
$data[64,9] = 15
$data[64,10] = 44000
$data[64,11] = 44000
$data[64,12] = 44000
$data[64,13] = '$/caja 18 kilos'
$data[64,14] = 'Provincia de Limarí'
$data[64,15] = 2444
$data[64,16] = 18
$data[64,17] = 'Hortaliza'
# target row 67 (was row 58)
$data[65,0] = 12
$data[65,1] = 'Mapocho Venta Directa de Santiago'
$data[65,2] = 'Metropolitana'
$data[65,3] = 44467
$data[65,4] = 13
$data[65,5] = 100112002
$data[65,6] = 'Pimiento'
$data[65,7] = 'Cuatro cascos verde'
$data[65,8] = 'Segunda'
$data[65,9] = 25
$data[65,10] = 42000
$data[65,11] = 42000
$data[65,12] = 42000
$data[65,13] = '$/caja 18 kilos'
$data[65,14] = 'Provincia de Limarí'
$data[65,15] = 2333
$data[65,16] = 18
$data[65,17] = 'Hortaliza'
# target row 68 (was row 59)
$data[66,0] = 12
$data[66,1] = 'Mapocho Venta Directa de Santiago'
$data[66,2] = 'Metropolitana'
$data[66,3] = 44467
$data[66,4] = 13
$data[66,5] = 100112002
$data[66,6] = 'Pimiento'
$data[66,7] = 'Morrón rojo'
$data[66,8] = 'Primera'
$data[66,9] = 10
$data[66,10] = 45000
$data[66,11] = 45000
$data[66,12] = 45000
$data[66,13] = '$/caja 18 kilos'
$data[66,14] = 'Provincia de Limarí'
$data[66,15] = 2500
$data[66,16] = 18
$data[66,17] = 'Hortaliza'
# target row 69 (was row 60)
$data[67,0] = 12
$data[67,1] = 'Mapocho Venta Directa de Santiago'
$data[67,2] = 'Metropolitana'
$data[67,3] = 44467
$data[67,4] = 13
$data[67,5] = 100112002
$data[67,6] = 'Pimiento'
$data[67,7] = 'Morrón rojo'
$data[67,8] = 'Segunda'
$data[67,9] = 15
$data[67,10] = 43000
$data[67,11] = 43000
$data[67,12] = 43000
$data[67,13] = '$/caja 18 kilos'
$data[67,14] = 'Provincia de Limarí'
$data[67,15] = 2389
$data[67,16] = 18
$data[67,17] = 'Hortaliza'
# target row 70 (was row 41)
$data[68,0] = 12
$data[68,1] = 'Mapocho Venta Directa de Santiago'
$data[68,2] = 'Metropolitana'
$data[68,3] = 44200
$data[68,4] = 13
$data[68,5] = 100112002
$data[68,6] = 'Pimiento'
$data[68,7] = 'Zafiro verde'
$data[68,8] = 'Primera'
$data[68,9] = 15
$data[68,10] = 16000
$data[68,11] = 16000
$data[68,12] = 16000
$data[68,13] = '$/caja 18 kilos'
$data[68,14] = 'Provincia de Limarí'
$data[68,15] = 889
$data[68,16] = 18
$data[68,17] = 'Hortaliza'
# target row 71 (was row 42)
$data[69,0] = 12
$data[69,1] = 'Mapocho Venta Directa de Santiago'
$data[69,2] = 'Metropolitana'
$data[69,3] = 44200
$data[69,4] = 13
$data[69,5] = 100112002
$data[69,6] = 'Pimiento'
$data[69,7] = 'Zafiro verde'
$data[69,8] = 'Segunda'
$data[69,9] = 10
$data[69,10] = 13000
$data[69,11] = 13000
$data[69,12] = 13000
$data[69,13] = '$/caja 18 kilos'
$data[69,14] = 'Provincia de Limarí'
$data[69,15] = 722
$data[69,16] = 18
$data[69,17] = 'Hortaliza'
# target row 72 (was row 83)
$data[70,0] = 12
$data[70,1] = 'Mapocho Venta Directa de Santiago'
$data[70,2] = 'Metropolitana'
$data[70,3] = 44544
$data[70,4] = 13
$data[70,5] = 100112002
$data[70,6] = 'Pimiento'
$data[70,7] = 'Morrón rojo'
$data[70,8] = 'Primera'
$data[70,9] = 10
$data[70,10] = 32000
$data[70,11] = 32000
$data[70,12] = 32000
$data[70,13] = '$/caja 18 kilos'
$data[70,14] = 'Provincia de Limarí'
$data[70,15] = 1778
$data[70,16] = 18
$data[70,17] = 'Hortaliza'
# target row 73 (was row 84)
$data[71,0] = 12
$data[71,1] = 'Mapocho Venta Directa de Santiago'
$data[71,2] = 'Metropolitana'
$data[71,3] = 44544
$data[71,4] = 13
$data[71,5] = 100112002
$data[71,6] = 'Pimiento'
$data[71,7] = 'Zafiro rojo'
$data[71,8] = 'Primera'
$data[71,9] = 15
$data[71,10] = 21000
$data[71,11] = 21000
$data[71,12] = 21000
$data[71,13] = '$/caja 18 kilos'
$data[71,14] = 'Provincia de Limarí'
$data[71,15] = 1167
$data[71,16] = 18
$data[71,17] = 'Hortaliza'
# target row 74 (was row 85)
$data[72,0] = 12
$data[72,1] = 'Mapocho Venta Directa de Santiago'
$data[72,2] = 'Metropolitana'
$data[72,3] = 44544
$data[72,4] = 13
$data[72,5] = 100112002
$data[72,6] = 'Pimiento'
$data[72,7] = 'Zafiro verde'
$data[72,8] = 'Primera'
$data[72,9] = 20
$data[72,10] = 15000
$data[72,11] = 15000
$data[72,12] = 15000
$data[72,13] = '$/caja 18 kilos'
$data[72,14] = 'Provincia de Limarí'
$data[72,15] = 833
$data[72,16] = 18
$data[72,17] = 'Hortaliza'
# target row 75 (was row 62)
$data[73,0] = 12
$data[73,1] = 'Mapocho Venta Directa de Santiago'
$data[73,2] = 'Metropolitana'
$data[73,3] = 44449
$data[73,4] = 13
$data[73,5] = 100112002
$data[73,6] = 'Pimiento'
$data[73,7] = 'Zafiro rojo'
$data[73,8] = 'Tercera'
$data[73,9] = 25
$data[73,10] = 34000
$data[73,11] = 34000
$data[73,12] = 34000
$data[73,13] = '$/caja 18 kilos'
$data[73,14] = 'Provincia de Limarí'
$data[73,15] = 1889
$data[73,16] = 18
$data[73,17] = 'Hortaliza'
# target row 76 (was row 63)
$data[74,0] = 12
$data[74,1] = 'Mapocho Venta Directa de Santiago'
$data[74,2] = 'Metropolitana'
$data[74,3] = 44449
$data[74,4] = 13
$data[74,5] = 100112002
$data[74,6] = 'Pimiento'
$data[74,7] = 'Zafiro verde'
$data[74,8] = 'Primera'
$data[74,9] = 15
$data[74,10] = 38000
$data[74,11] = 38000
$data[74,12] = 38000
$data[74,13] = '$/caja 18 kilos'
$data[74,14] = 'Provincia de Limarí'
$data[74,15] = 2111
$data[74,16] = 18
$data[74,17] = 'Hortaliza'
# target row 77 (was row 64)
$data[75,0] = 12
$data[75,1] = 'Mapocho Venta Directa de Santiago'
$data[75,2] = 'Metropolitana'
$data[75,3] = 44449
$data[75,4] = 13
$data[75,5] = 100112002
$data[75,6] = 'Pimiento'
$data[75,7] = 'Zafiro verde'
$data[75,8] = 'Segunda'
$data[75,9] = 20
$data[75,10] = 36000
$data[75,11] = 36000
$data[75,12] = 36000
$data[75,13] = '$/caja 18 kilos'
$data[75,14] = 'Provincia de Limarí'
$data[75,15] = 2000
$data[75,16] = 18
$data[75,17] = 'Hortaliza'
# target row 78 (was row 65)
$data[76,0] = 12
$data[76,1] = 'Mapocho Venta Directa de Santiago'
$data[76,2] = 'Metropolitana'
$data[76,3] = 44449
$data[76,4] = 13
$data[76,5] = 100112002
$data[76,6] = 'Pimiento'
$data[76,7] = 'Zafiro verde'
$data[76,8] = 'Tercera'
$data[76,9] = 25
$data[76,10] = 34000
$data[76,11] = 34000
$data[76,12] = 34000
$data[76,13] = '$/caja 18 kilos'
$data[76,14] = 'Provincia de Limarí'
$data[76,15] = 1889
$data[76,16] = 18
$data[76,17] = 'Hortaliza'
# target row 79 (was row 18)
$data[77,0] = 12
$data[77,1] = 'Mapocho Venta Directa de Santiago'
$data[77,2] = 'Metropolitana'
$data[77,3] = 44165
$data[77,4] = 13
$data[77,5] = 100112002
$data[77,6] = 'Pimiento'
$data[77,7] = 'Zafiro rojo'
$data[77,8] = 'Primera'
$data[77,9] = 125
$data[77,10] = 30000
$data[77,11] = 30000
$data[77,12] = 30000
$data[77,13] = '$/caja 15 kilos'
$data[77,14] = 'Región de Arica y Parinacota'
$data[77,15] = 2000
$data[77,16] = 15
$data[77,17] = 'Hortaliza'
# target row 80 (was row 19)
$data[78,0] = 12
$data[78,1] = 'Mapocho Venta Directa de Santiago'
$data[78,2] = 'Metropolitana'
$data[78,3] = 44165
$data[78,4] = 13
$data[78,5] = 100112002
$data[78,6] = 'Pimiento'
$data[78,7] = 'Zafiro verde'
$data[78,8] = 'Primera'
$data[78,9] = 60
$data[78,10] = 28000
$data[78,11] = 28000
$data[78,12] = 28000
$data[78,13] = '$/caja 18 kilos'
$data[78,14] = 'Provincia de Quillota'
$data[78,15] = 1556
$data[78,16] = 18
$data[78,17] = 'Hortaliza'
# target row 81 (was row 20)
$data[79,0] = 12
$data[79,1] = 'Mapocho Venta Directa de Santiago'
$data[79,2] = 'Metropolitana'
$data[79,3] = 44165
$data[79,4] = 13
$data[79,5] = 100112002
$data[79,6] = 'Pimiento'
$data[79,7] = 'Zafiro verde'
$data[79,8] = 'Segunda'
$data[79,9] = 40
$data[79,10] = 25000
$data[79,11] = 25000
$data[79,12] = 25000
$data[79,13] = '$/caja 18 kilos'
$data[79,14] = 'Provincia de Quillota'
$data[79,15] = 1389
$data[79,16] = 18
$data[79,17] = 'Hortaliza'
# target row 82 (was row 55)
$data[80,0] = 12
$data[80,1] = 'Mapocho Venta Directa de Santiago'
$data[80,2] = 'Metropolitana'
$data[80,3] = 44350
$data[80,4] = 13
$data[80,5] = 100112002
$data[80,6] = 'Pimiento'
$data[80,7] = 'Zafiro rojo'
$data[80,8] = 'Primera'
$data[80,9] = 15
$data[80,10] = 28000
$data[80,11] = 28000
$data[80,12] = 28000
$data[80,13] = '$/caja 15 kilos'
$data[80,14] = 'Región de Arica y Parinacota'
$data[80,15] = 1867
$data[80,16] = 15
$data[80,17] = 'Hortaliza'
# target row 83 (was row 56)
$data[81,0] = 12
$data[81,1] = 'Mapocho Venta Directa de Santiago'
$data[81,2] = 'Metropolitana'
$data[81,3] = 44350
$data[81,4] = 13
$data[81,5] = 100112002
$data[81,6] = 'Pimiento'
$data[81,7] = 'Zafiro verde'
$data[81,8] = 'Primera'
$data[81,9] = 35
$data[81,10] = 14000
$data[81,11] = 14000
$data[81,12] = 14000
$data[81,13] = '$/caja 15 kilos'
$data[81,14] = 'Región de Arica y Parinacota'
$data[81,15] = 933
$data[81,16] = 15
$data[81,17] = 'Hortaliza'
# target row 84 (was row 26)
$data[82,0] = 12
$data[82,1] = 'Mapocho Venta Directa de Santiago'
$data[82,2] = 'Metropolitana'
$data[82,3] = 44305
$data[82,4] = 13
$data[82,5] = 100112002
$data[82,6] = 'Pimiento'
$data[82,7] = 'Zafiro rojo'
$data[82,8] = 'Primera'
$data[82,9] = 35
$data[82,10] = 17000
$data[82,11] = 18000
$data[82,12] = 17571
$data[82,13] = '$/caja 18 kilos'
$data[82,14] = 'Provincia de Limarí'
$data[82,15] = 976
$data[82,16] = 18
$data[82,17] = 'Hortaliza'
# target row 85 (was row 27)
$data[83,0] = 12
$data[83,1] = 'Mapocho Venta Directa de Santiago'
$data[83,2] = 'Metropolitana'
$data[83,3] = 44305
$data[83,4] = 13
$data[83,5] = 100112002
$data[83,6] = 'Pimiento'
$data[83,7] = 'Zafiro verde'
$data[83,8] = 'Primera'
$data[83,9] = 60
$data[83,10] = 14000
$data[83,11] = 14000
$data[83,12] = 14000
$data[83,13] = '$/caja 18 kilos'
$data[83,14] = 'Provincia de Limarí'
$data[83,15] = 778
$data[83,16] = 18
$data[83,17] = 'Hortaliza'
# target row 86 (was row 51)
$data[84,0] = 12
$data[84,1] = 'Mapocho Venta Directa de Santiago'
$data[84,2] = 'Metropolitana'
$data[84,3] = 44273
$data[84,4] = 13
$data[84,5] = 100112002
$data[84,6] = 'Pimiento'
$data[84,7] = 'Zafiro rojo'
$data[84,8] = 'Primera'
$data[84,9] = 45
$data[84,10] = 17000
$data[84,11] = 17000
$data[84,12] = 17000
$data[84,13] = '$/caja 18 kilos'
$data[84,14] = 'Provincia de Limarí'
$data[84,15] = 944
$data[84,16] = 18
$data[84,17] = 'Hortaliza'
# target row 87 (was row 52)
$data[85,0] = 12
$data[85,1] = 'Mapocho Venta Directa de Santiago'
$data[85,2] = 'Metropolitana'
$data[85,3] = 44273
$data[85,4] = 13
$data[85,5] = 100112002
$data[85,6] = 'Pimiento'
$data[85,7] = 'Zafiro rojo'
$data[85,8] = 'Segunda'
$data[85,9] = 30
$data[85,10] = 15000
$data[85,11] = 15000
$data[85,12] = 15000
$data[85,13] = '$/caja 18 kilos'
$data[85,14] = 'Provincia de Limarí'
$data[85,15] = 833
$data[85,16] = 18
$data[85,17] = 'Hortaliza'
# target row 88 (was row 53)
$data[86,0] = 12
$data[86,1] = 'Mapocho Venta Directa de Santiago'
$data[86,2] = 'Metropolitana'
$data[86,3] = 44273
$data[86,4] = 13
$data[86,5] = 100112002
$data[86,6] = 'Pimiento'
$data[86,7] = 'Zafiro rojo'
$data[86,8] = 'Tercera'
$data[86,9] = 25
$data[86,10] = 13000
$data[86,11] = 13000
$data[86,12] = 13000
$data[86,13] = '$/caja 18 kilos'
$data[86,14] = 'Provincia de Limarí'
$data[86,15] = 722
$data[86,16] = 18
$data[86,17] = 'Hortaliza'

$ws.Range("A2:R88").Value = $data

Write-Host "Rewrote A2:R88 with the refreshed weekly row ordering."
